$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert the new BOM line item (DWM1000 module) as row 15, matching the
# formatting of the row above it (row 14).
$ws.Range("A14:I14").Copy()
$ws.Range("A15:I15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A15").Value = 9
$ws.Range("B15").Value = "DWM1000"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "Decawave"
$ws.Range("E15").Value = "DWM1000"
$ws.Range("F15").Value = "IC RF-Module UWB Transceiver"
$ws.Range("G15").Value = ""
$ws.Range("H15").Value = "SMD"
$ws.Range("I15").Value = ""

# Keep view/selection similar to the saved file
$ws.Range("C16").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
